# "Completed the driver part"
# - Orders 3 and 4 (rows 4 and 5) have finished cooking and were handed off
#   to a driver, so their order status moves from "Food Prepared" to
#   "Completed" and they now have a Driver ID assigned (8).
# - The active selection is left on F9, matching where the author's cursor
#   ended up after making the edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order status: Food Prepared -> Completed
$ws.Range("E4").Value = "Completed"
$ws.Range("E5").Value = "Completed"

# Driver ID assigned now that the order is ready for delivery
$ws.Range("I4").Value = 8
$ws.Range("I5").Value = 8

# Leave the selection where the author left it
$ws.Range("F9").Select()
